$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 14
$ws.Range("B3").Value = 106471685
$ws.Range("C3").Value = 'G'
$ws.Range("D3").Value = 'A'
$ws.Range("F3").Value = $false
$ws.Range("G3").Value = 'Mapping problems a plenty, variant reads all have same pattern of mismatches, lots of raw variant calls made by Caveman in this region that get filtered because of strand bias'
$ws.Range("H3").Value = 'Zoom out to see lots of mismapping reads and lots of filtered and unfiltered variant calls if VCF is loaded as a track'

# Row 4
$ws.Range("A4").Value = 15
$ws.Range("B4").Value = 91537757
$ws.Range("C4").Value = 'C'
$ws.Range("D4").Value = 'G'
$ws.Range("F4").Value = $false
$ws.Range("G4").Value = 'Filtered by cgpCaVEManPostProcessor, MNP filter (tumour AF - normal AF < 0.2), 4 variant supporting reads in normal but 3 are very low base quality'
$ws.Range("H4").ClearContents()

# Row 5
$ws.Range("A5").Value = 22
$ws.Range("B5").Value = 38120429
$ws.Range("C5").Value = 'T'
$ws.Range("D5").Value = 'C'
$ws.Range("F5").Value = '?'
$ws.Range("G5").Value = 'Filtered by cgpCaVEManPostProcessor, MQ filter (mean mapping quality < 21), reads with zero mapping quality in vicinity which could map equally well to another place in the genome'
$ws.Range("H5").Value = 'Zoom out to view zero mapping quality reads in vicinity; sort by mapping quality at the variant position'

# Row 6
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = 9983580
$ws.Range("C6").Value = 'T'
$ws.Range("D6").Value = 'C'
$ws.Range("F6").Value = $false
$ws.Range("G6").Value = 'Filtered by cgpCaVEManPostProcessor, RP filter (no mutant alleles found in first 2/3 of read and coverage < 8), misalignment of reads that just span indel, would have been resolved by a caller that performs local assembly around variant regions'
$ws.Range("H6").ClearContents()

# Row 7
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = 9633452
$ws.Range("C7").Value = 'G'
$ws.Range("D7").Value = 'A'
$ws.Range("F7").Value = $true
$ws.Range("G7").Value = 'Low allele fraction (13%, 4/31), missense and predicted to be damaging, called by MuTect2 but not Caveman'
$ws.Range("H7").ClearContents()

# Row 8
$ws.Range("A8").Value = 3
$ws.Range("B8").Value = 177090982
$ws.Range("C8").Value = 'A'
$ws.Range("D8").Value = 'G'
$ws.Range("F8").Value = '?'
$ws.Range("G8").Value = 'Low depth in the normal is the main concern here - is this really somatic or could it be germline? Probability of having 6 of 6 reads with reference allele if this was a heterozygous germline variant is 0.016 but note that several germline positions may only have been sequenced to low depth in the normal (dangers of not sequencing the normal to sufficient depth)'
$ws.Range("H8").ClearContents()

# Row 9
$ws.Range("A9").Value = 15
$ws.Range("B9").Value = 68175783
$ws.Range("C9").Value = 'A'
$ws.Range("D9").Value = 'G'
$ws.Range("F9").Value = $false
$ws.Range("G9").Value = 'Filtered because of strand bias, 10 variant supporting reads all on forward strand'
$ws.Range("H9").Value = 'Sort alignments by base, colour by read strand'

# Row 10
$ws.Range("A10").Value = 3
$ws.Range("B10").Value = 140281877
$ws.Range("C10").Value = 'A'
$ws.Range("D10").Value = 'C'
$ws.Range("F10").Value = '?'
$ws.Range("G10").Value = 'Filtered by cgpCaVEManPostProcessor, DTH filter (< 1/3 mutant alleles with base quality >= 25), 2 variant supporting reads in normal but with very low base qualities'
$ws.Range("H10").ClearContents()

# Row heights
$ws.Rows(5).RowHeight = 45
$ws.Rows(6).RowHeight = 45
$ws.Rows(7).RowHeight = 30
$ws.Rows(8).RowHeight = 75
$ws.Rows(9).RowHeight = 15
$ws.Rows(10).RowHeight = 30

# Update selection to F11
$ws.Range("F11").Select()
